$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "306.70" or
# "42.102.06" keep their exact formatting instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.102.06"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "2.250.41"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "306.70"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").Value = "96.91"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("E7").Value = "  -1.24%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("D10").Value = "35.09"
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("E11").Value = "  +2.54%  "
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("D13").Value = "6.76"
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("D14").Value = "2.598.85"
$ws.Range("E14").Value = "  -1.01%  "
$ws.Range("D15").Value = "14.45"
$ws.Range("E15").Value = "  +0.56%  "
$ws.Range("D16").Value = "2.259.73"
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("D17").Value = "0.780"
$ws.Range("E17").Value = "  -1.83%  "
$ws.Range("D18").Value = "42.056.22"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").Value = "12.18"
$ws.Range("E19").Value = "  -2.69%  "
$ws.Range("D20").Value = "0.0₃0903"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").Value = "5.93"
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("D22").Value = "67.14"
$ws.Range("E22").Value = "  -0.65%  "
$ws.Range("D23").Value = "236.10"
$ws.Range("E23").Value = "  -1.81%  "
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("B27").Value = "InjectiveProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D27").Value = "38.01"
$ws.Range("E27").Value = "  +1.82%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "23.34"
$ws.Range("E28").Value = "  -1.99%  "
$ws.Range("D29").Value = "9.49"
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("D31").Value = "167.05"
$ws.Range("E31").Value = "  +4.43%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("E33").Value = "  -1.69%  "
$ws.Range("D34").Value = "17.56"
$ws.Range("E34").Value = "  +2.87%  "
$ws.Range("E35").Value = "  -3.59%  "
$ws.Range("D36").Value = "0.0720"
$ws.Range("E36").Value = "  -2.96%  "
$ws.Range("E37").Value = "  +2.04%  "
$ws.Range("D38").Value = "0.115"
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("D39").Value = "0.103"
$ws.Range("E39").Value = "  -2.64%  "
$ws.Range("E40").Value = "  -2.17%  "
$ws.Range("D41").Value = "4.08"
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").Value = "1.940.70"
$ws.Range("E42").Value = "  -2.82%  "
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("D44").Value = "18.61"
$ws.Range("E44").Value = "  -0.81%  "
$ws.Range("D45").Value = "2.18"
$ws.Range("E45").Value = "  -10.42%  "
$ws.Range("D46").Value = "2.90"
$ws.Range("E46").Value = "  -1.52%  "
$ws.Range("D47").Value = "9.66"
$ws.Range("E47").Value = "  -3.19%  "
$ws.Range("D48").Value = "54.15"
$ws.Range("E48").Value = "  +1.71%  "
$ws.Range("D49").Value = "2.468.86"
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("D50").Value = "71.15"
$ws.Range("E50").Value = "  -1.44%  "
$ws.Range("D51").Value = "91.17"
$ws.Range("E51").Value = "  -0.42%  "
